# Applies the row rotation described in the commit:
#   - Row 7  now holds the record previously found in Row 9
#   - Row 9  now holds the record previously found in Row 11
#   - Row 10 now holds the record previously found in Row 7 (original)
#   - Row 11 now holds the record previously found in Row 10 (original)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7 : Knärot (112212882) -> Ullticka (112213279) ----
$ws.Range("A7").Value2 = 112213279
$ws.Range("B7").Value2 = 89553
$ws.Range("D7").Value2 = "NT"
$ws.Range("E7").Value2 = 1202
$ws.Range("F7").Value2 = "Ullticka"
$ws.Range("G7").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("P7").Value2 = "Nordvallen (Nordvallen), Jmt"
$ws.Range("Q7").Value2 = 490080
$ws.Range("R7").Value2 = 6948907

# ---- Row 9 : Ullticka (112213279) -> Knärot (112212836) ----
$ws.Range("A9").Value2 = 112212836
$ws.Range("B9").Value2 = 96735
$ws.Range("D9").Value2 = "VU"
$ws.Range("E9").Value2 = 220787
$ws.Range("F9").Value2 = "Knärot"
$ws.Range("G9").Value2 = "Goodyera repens"
$ws.Range("H9").Value2 = "(L.) R. Br."
$ws.Range("I9").Value2 = 25
$ws.Range("J9").Value2 = "plantor/tuvor"
$ws.Range("P9").Value2 = "Stugunäset (Stugunäset), Jmt"
$ws.Range("Q9").Value2 = 490078
$ws.Range("R9").Value2 = 6948752

# ---- Row 10 : 112212105 -> 112212882 ----
$ws.Range("A10").Value2 = 112212882
$ws.Range("J10").Value2 = "plantor/tuvor"
$ws.Range("P10").Value2 = "Kälen (Kälen), Jmt"
$ws.Range("Q10").Value2 = 490109
$ws.Range("R10").Value2 = 6948768
$ws.Range("Z10").Value2 = "12:39"
$ws.Range("AB10").Value2 = "12:39"

# ---- Row 11 : 112212836 -> 112212105 ----
$ws.Range("A11").Value2 = 112212105
$ws.Range("I11").Value2 = 3
$ws.Range("J11").ClearContents()
$ws.Range("P11").Value2 = "Nordvallen (Nordvallen), Jmt"
$ws.Range("Q11").Value2 = 490018
$ws.Range("R11").Value2 = 6948882
$ws.Range("Z11").Value2 = "11:58"
$ws.Range("AB11").Value2 = "11:58"
